$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

$ws.Range("P5").Select()
$win.ScrollColumn = 16
$win.ScrollRow = 5
$ws.Range("A27").Select()
